# Bonus zones for 15/30
# Adds a second target species column (TARGET_2) and 16 new "bonus" sampling
# zone rows (stations 26-41), while re-numbering the existing CLUSTER ids to
# make room for the new clusters.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add TARGET_2 column
$ws.Cells.Item(1, 4).Value = "TARGET_2"

# Update CLUSTER numbers (column B) and station/target data for rows 2-25
$ws.Cells.Item(2, 1).Value = "N-23"
$ws.Cells.Item(2, 2).Value = 3
$ws.Cells.Item(2, 3).Value = "snow crab"

$ws.Cells.Item(3, 1).Value = "N-22"
$ws.Cells.Item(3, 2).Value = 3
$ws.Cells.Item(3, 3).Value = "snow crab"

$ws.Cells.Item(4, 1).Value = "N-21"
$ws.Cells.Item(4, 2).Value = 3
$ws.Cells.Item(4, 3).Value = "snow crab"

$ws.Cells.Item(5, 1).Value = "N-20"
$ws.Cells.Item(5, 2).Value = 3
$ws.Cells.Item(5, 3).Value = "snow crab"

$ws.Cells.Item(6, 1).Value = "K-22"
$ws.Cells.Item(6, 2).Value = 4
$ws.Cells.Item(6, 3).Value = "snow crab"

$ws.Cells.Item(7, 1).Value = "K-21"
$ws.Cells.Item(7, 2).Value = 4
$ws.Cells.Item(7, 3).Value = "snow crab"

$ws.Cells.Item(8, 1).Value = "K-20"
$ws.Cells.Item(8, 2).Value = 4
$ws.Cells.Item(8, 3).Value = "snow crab"

$ws.Cells.Item(9, 1).Value = "K-19"
$ws.Cells.Item(9, 2).Value = 4
$ws.Cells.Item(9, 3).Value = "snow crab"

$ws.Cells.Item(10, 1).Value = "I-22"
$ws.Cells.Item(10, 2).Value = 5
$ws.Cells.Item(10, 3).Value = "snow crab"

$ws.Cells.Item(11, 1).Value = "I-21"
$ws.Cells.Item(11, 2).Value = 5
$ws.Cells.Item(11, 3).Value = "snow crab"

$ws.Cells.Item(12, 1).Value = "I-20"
$ws.Cells.Item(12, 2).Value = 5
$ws.Cells.Item(12, 3).Value = "snow crab"

$ws.Cells.Item(13, 1).Value = "H-20"
$ws.Cells.Item(13, 2).Value = 5
$ws.Cells.Item(13, 3).Value = "snow crab"

$ws.Cells.Item(14, 1).Value = "G-11"
$ws.Cells.Item(14, 2).Value = 10
$ws.Cells.Item(14, 3).Value = "red king crab"

$ws.Cells.Item(15, 1).Value = "F-11"
$ws.Cells.Item(15, 2).Value = 10
$ws.Cells.Item(15, 3).Value = "red king crab"

$ws.Cells.Item(16, 1).Value = "E-11"
$ws.Cells.Item(16, 2).Value = 10
$ws.Cells.Item(16, 3).Value = "red king crab"

$ws.Cells.Item(17, 1).Value = "E-12"
$ws.Cells.Item(17, 2).Value = 10
$ws.Cells.Item(17, 3).Value = "red king crab"

$ws.Cells.Item(18, 1).Value = "G-22"
$ws.Cells.Item(18, 2).Value = 6
$ws.Cells.Item(18, 3).Value = "Tanner crab"
$ws.Cells.Item(18, 4).Value = "Alaska skate"

$ws.Cells.Item(19, 1).Value = "G-21"
$ws.Cells.Item(19, 2).Value = 6
$ws.Cells.Item(19, 3).Value = "Tanner crab"
$ws.Cells.Item(19, 4).Value = "Alaska skate"

$ws.Cells.Item(20, 1).Value = "G-20"
$ws.Cells.Item(20, 2).Value = 6
$ws.Cells.Item(20, 3).Value = "Tanner crab"
$ws.Cells.Item(20, 4).Value = "Alaska skate"

$ws.Cells.Item(21, 1).Value = "G-19"
$ws.Cells.Item(21, 2).Value = 6
$ws.Cells.Item(21, 3).Value = "Tanner crab"
$ws.Cells.Item(21, 4).Value = "Alaska skate"

$ws.Cells.Item(22, 1).Value = "D-04"
$ws.Cells.Item(22, 2).Value = 7
$ws.Cells.Item(22, 3).Value = "arrowtooth flounder"

$ws.Cells.Item(23, 1).Value = "D-05"
$ws.Cells.Item(23, 2).Value = 7
$ws.Cells.Item(23, 3).Value = "arrowtooth flounder"

$ws.Cells.Item(24, 1).Value = "C-05"
$ws.Cells.Item(24, 2).Value = 7
$ws.Cells.Item(24, 3).Value = "arrowtooth flounder"

$ws.Cells.Item(25, 1).Value = "B-05"
$ws.Cells.Item(25, 2).Value = 7
$ws.Cells.Item(25, 3).Value = "arrowtooth flounder"

# New bonus sampling zone rows 26-41.
# Column A (station IDs) is filled in first for the whole new "Crab" block,
# then column C (target) is filled in afterwards, matching how the
# new rows were originally populated (station list pasted first, target
# label filled down after) - rows 38-41 (Tanner crab group) repeat the
# same pattern as a second, later pass.
$ws.Cells.Item(26, 1).Value = "U-29"
$ws.Cells.Item(27, 1).Value = "U-28"
$ws.Cells.Item(28, 1).Value = "U-27"
$ws.Cells.Item(29, 1).Value = "U-26"
$ws.Cells.Item(30, 1).Value = "J-01"
$ws.Cells.Item(31, 1).Value = "I-01"
$ws.Cells.Item(32, 1).Value = "H-01"
$ws.Cells.Item(33, 1).Value = "G-01"
$ws.Cells.Item(34, 1).Value = "H-09"
$ws.Cells.Item(35, 1).Value = "H-08"
$ws.Cells.Item(36, 1).Value = "H-07"
$ws.Cells.Item(37, 1).Value = "H-06"

$ws.Cells.Item(26, 3).Value = "Crab"
$ws.Cells.Item(27, 3).Value = "Crab"
$ws.Cells.Item(28, 3).Value = "Crab"
$ws.Cells.Item(29, 3).Value = "Crab"
$ws.Cells.Item(30, 3).Value = "Crab"
$ws.Cells.Item(31, 3).Value = "Crab"
$ws.Cells.Item(32, 3).Value = "Crab"
$ws.Cells.Item(33, 3).Value = "Crab"
$ws.Cells.Item(34, 3).Value = "Crab"
$ws.Cells.Item(35, 3).Value = "Crab"
$ws.Cells.Item(36, 3).Value = "Crab"
$ws.Cells.Item(37, 3).Value = "Crab"

$ws.Cells.Item(38, 1).Value = "M-30"
$ws.Cells.Item(39, 1).Value = "M-29"
$ws.Cells.Item(40, 1).Value = "M-28"
$ws.Cells.Item(41, 1).Value = "M-27"

$ws.Cells.Item(38, 3).Value = "Tanner crab"
$ws.Cells.Item(39, 3).Value = "Tanner crab"
$ws.Cells.Item(40, 3).Value = "Tanner crab"
$ws.Cells.Item(41, 3).Value = "Tanner crab"

# Column B (CLUSTER numbers) for the new rows
$ws.Cells.Item(26, 2).Value = 1
$ws.Cells.Item(27, 2).Value = 1
$ws.Cells.Item(28, 2).Value = 1
$ws.Cells.Item(29, 2).Value = 1
$ws.Cells.Item(30, 2).Value = 8
$ws.Cells.Item(31, 2).Value = 8
$ws.Cells.Item(32, 2).Value = 8
$ws.Cells.Item(33, 2).Value = 8
$ws.Cells.Item(34, 2).Value = 9
$ws.Cells.Item(35, 2).Value = 9
$ws.Cells.Item(36, 2).Value = 9
$ws.Cells.Item(37, 2).Value = 9
$ws.Cells.Item(38, 2).Value = 2
$ws.Cells.Item(39, 2).Value = 2
$ws.Cells.Item(40, 2).Value = 2
$ws.Cells.Item(41, 2).Value = 2

# Restore the view/selection state (scrolled up a bit, new active cell)
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("C23").Select()
